$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 (Tiempo_Mínimo, Tiempo_Máximo, Tiempo_Promedio)
$ws.Range("E2").Value = 0.0001629
$ws.Range("F2").Value = 0.01434411
$ws.Range("G2").Value = 0.00026464986027190335

# Row 3
$ws.Range("E3").Value = 0.00179442
$ws.Range("F3").Value = 0.01600605
$ws.Range("G3").Value = 0.0022136484611973394

# Row 4
$ws.Range("E4").Value = 0.00819099
$ws.Range("F4").Value = 0.01660491
$ws.Range("G4").Value = 0.009364566741573033
